# Apply updated crypto price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.659.70"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "2.511.91"
$ws.Range("E3").Value = "  -1.62%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "316.35"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.17%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "95.75"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "

$ws.Range("E7").Value = "  +1.69%  "

$ws.Range("E8").Value = "  -0.04%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.538"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.69%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "36.36"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.16%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0813"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "7.73"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.39%  "

$ws.Range("E13").Value = "  -3.35%  "

$ws.Range("D14").Value = "2.904.96"
$ws.Range("E14").Value = "  -1.40%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.54"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +4.66%  "

$ws.Range("D16").Value = "2.530.97"
$ws.Range("E16").Value = "  -0.31%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.864"
$c.Style = "Normal"

$ws.Range("D18").Value = "42.668.48"
$ws.Range("E18").Value = "  -0.39%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.88"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.69%  "

$ws.Range("D20").Value = "0.0₃0974"
$ws.Range("E20").Value = "  -1.13%  "

$ws.Range("E21").Value = "  -0.13%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "71.49"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "253.19"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.99"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.91%  "

$ws.Range("E25").Value = "  -2.39%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "27.02"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.54%  "

$ws.Range("E27").Value = "  -0.02%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +13.13%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "10.13"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.45%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "38.05"
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "5.92"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.17%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "156.20"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.30%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "19.63"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +5.59%  "

$ws.Range("E34").Value = "  +1.42%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.08"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.43%  "

$ws.Range("E36").Value = "  -2.53%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.62"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.74%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.113"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.08%  "

$ws.Range("E39").Value = "  +1.00%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "24.20"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -9.13%  "

$ws.Range("E41").Value = "  -0.04%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.38"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("E43").Value = "  -2.15%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0303"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D46").Value = "2.036.99"
$ws.Range("E46").Value = "  -1.23%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "84.51"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.43%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "8.96"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.35%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "74.71"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "

$ws.Range("D50").Value = "2.761.29"
$ws.Range("E50").Value = "  -1.48%  "

$ws.Range("E51").Value = "  +0.20%  "
